# Re-sync per-batch stock rows (B/E/F/G columns: batch no., rate, qty, value)
# for several items whose batch lines were reordered, and refresh the
# downstream sub-totals / grand totals that depend on them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B149").Value = 63902
$ws.Range("E149").Value = 34.04
$ws.Range("F149").Value = 2
$ws.Range("G149").Value = 64.04000000000001

$ws.Range("B150").Value = 48654
$ws.Range("E150").Value = 38.26
$ws.Range("F150").Value = -1
$ws.Range("G150").Value = -32.02

$ws.Range("B279").Value = 64973
$ws.Range("E279").Value = 35.4
$ws.Range("F279").Value = 145
$ws.Range("G279").Value = 4828.5

$ws.Range("B280").Value = 48706
$ws.Range("E280").Value = 39.8
$ws.Range("F280").Value = -144
$ws.Range("G280").Value = -4795.2

$ws.Range("B313").Value = 62997
$ws.Range("F313").Value = 0
$ws.Range("G313").Value = 0

$ws.Range("B314").Value = 57854
$ws.Range("F314").Value = 2
$ws.Range("G314").Value = 611.6799999999999

$ws.Range("B316").Value = 61610
$ws.Range("E316").Value = 122.71
$ws.Range("F316").Value = -58
$ws.Range("G316").Value = -5957.18

$ws.Range("B317").Value = 57077
$ws.Range("D317").Value = 93.08
$ws.Range("E317").Value = 111.2
$ws.Range("F317").Value = 1
$ws.Range("G317").Value = 93.08

$ws.Range("B318").Value = 63565
$ws.Range("D318").Value = 102.71
$ws.Range("E318").Value = 109.19
$ws.Range("F318").Value = 60
$ws.Range("G318").Value = 6162.6

$ws.Range("B350").Value = 57802
$ws.Range("E350").Value = 162.71
$ws.Range("F350").Value = -79
$ws.Range("G350").Value = -11334.92

$ws.Range("B351").Value = 63531
$ws.Range("E351").Value = 152.53
$ws.Range("F351").Value = 80
$ws.Range("G351").Value = 11478.4

$ws.Range("B352").Value = 63571
$ws.Range("F352").Value = 18
$ws.Range("G352").Value = 2582.64

$ws.Range("B372").Value = 63652
$ws.Range("E372").Value = 55.42
$ws.Range("F372").Value = 204
$ws.Range("G372").Value = 10634.52

$ws.Range("B373").Value = 57885
$ws.Range("E373").Value = 62.28
$ws.Range("F373").Value = 4
$ws.Range("G373").Value = 208.52

$ws.Range("B382").Value = 60325
$ws.Range("E382").Value = 151.57
$ws.Range("F382").Value = -102
$ws.Range("G382").Value = -12939.72

$ws.Range("B383").Value = 63560
$ws.Range("E383").Value = 134.87
$ws.Range("F383").Value = 31
$ws.Range("G383").Value = 3932.66

$ws.Range("B400").Value = 62933
$ws.Range("F400").Value = 138
$ws.Range("G400").Value = 8159.94

$ws.Range("B401").Value = 57835
$ws.Range("F401").Value = 1
$ws.Range("G401").Value = 59.13

$ws.Range("B421").Value = 63008
$ws.Range("F421").Value = 449
$ws.Range("G421").Value = 67875.33

$ws.Range("B422").Value = 57857
$ws.Range("F422").Value = 3
$ws.Range("G422").Value = 453.51

$ws.Range("B536").Value = 47097
$ws.Range("D536").Value = 112.28
$ws.Range("E536").Value = 134.16
$ws.Range("F536").Value = 15
$ws.Range("G536").Value = 1684.2

$ws.Range("B537").Value = 58047
$ws.Range("D537").Value = 105.54
$ws.Range("E537").Value = 126.1
$ws.Range("F537").Value = 46
$ws.Range("G537").Value = 4854.84

$ws.Range("B583").Value = 53263
$ws.Range("E583").Value = 15.29
$ws.Range("F583").Value = -309
$ws.Range("G583").Value = -3958.29

$ws.Range("B584").Value = 65066
$ws.Range("E584").Value = 13.61
$ws.Range("F584").Value = 257
$ws.Range("G584").Value = 3292.17

$ws.Range("B586").Value = 64915
$ws.Range("E586").Value = 20.98
$ws.Range("F586").Value = 14
$ws.Range("G586").Value = 276.22

$ws.Range("B587").Value = 45695
$ws.Range("E587").Value = 23.58
$ws.Range("F587").Value = -36
$ws.Range("G587").Value = -710.28

$ws.Range("B590").Value = 64922
$ws.Range("E590").Value = 20.98
$ws.Range("F590").Value = 170
$ws.Range("G590").Value = 3354.1

$ws.Range("B591").Value = 45706
$ws.Range("E591").Value = 23.58
$ws.Range("F591").Value = -202
$ws.Range("G591").Value = -3985.46

$ws.Range("B593").Value = 45718
$ws.Range("E593").Value = 19.38
$ws.Range("F593").Value = -294
$ws.Range("G593").Value = -4768.68

$ws.Range("B594").Value = 64927
$ws.Range("E594").Value = 17.26
$ws.Range("F594").Value = 281
$ws.Range("G594").Value = 4557.82

$ws.Range("B599").Value = 45709
$ws.Range("E599").Value = 15.69
$ws.Range("F599").Value = -300
$ws.Range("G599").Value = -3945

$ws.Range("B600").Value = 64925
$ws.Range("E600").Value = 13.97
$ws.Range("F600").Value = 274
$ws.Range("G600").Value = 3603.1

$ws.Range("B601").Value = 64919
$ws.Range("E601").Value = 27.97
$ws.Range("F601").Value = 189
$ws.Range("G601").Value = 4970.7

$ws.Range("B602").Value = 45702
$ws.Range("E602").Value = 31.43
$ws.Range("F602").Value = -215
$ws.Range("G602").Value = -5654.5

$ws.Range("F629").Value = 321
$ws.Range("G629").Value = 7784.25

$ws.Range("B651").Value = 50066.12

$ws.Range("F872").Value = 94
$ws.Range("G872").Value = 23451.12

$ws.Range("F898").Value = 212
$ws.Range("G898").Value = 51261.6

$ws.Range("B900").Value = 202582.81

$ws.Range("B962").Value = 4917871.19

$ws.Range("B963").Value = 4917871.19
